$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last calibration block (rows 32:37, date 2020-11-01 / serial 44136) is
# being repeated for a new calibration run on 2021-02-09 (serial 44236).
# Insert a copy of that block right after the existing data (new rows 38:43)
# so the duplicated rows keep the exact same styling (date format on A,
# the "touched" formatting on B/C, etc.) that a copy/insert would carry in
# real Excel.
$ws.Rows("32:37").Copy()
$ws.Rows("38:43").Insert(-4121)  # xlShiftDown

# The newly inserted rows came in with the old calibration's start/end-time
# values (columns G/H) and date (column A) still attached - clear those so
# only date/plot/type/start_min/end_min carry over, matching the new,
# not-yet-timed calibration entries.
$ws.Range("G38:H43").Clear()
$ws.Range("A38:A43").Value = 44236

# The last row of the new block (the second plot's "amb" reading) used
# different start/end minutes than the block it was copied from - correct
# them here.
$ws.Cells.Item(43, 4).Value = 35
$ws.Cells.Item(43, 5).Value = 38

# Restore the user's current selection to where they left off editing.
$ws.Range("G37").Select()

# Calibrate the print/page setup for this sheet (portrait, A4).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
